$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: new header "BOUNDARY" merged across AS1:AZ1, styled like the other header blocks ---
# Merge first, then apply the header format/value so the whole merged area
# ends up with one uniform style (matching the other header blocks) instead
# of per-cell merge-boundary borders.
$ws.Range("AS1:AZ1").Merge()
$ws.Range("AK1:AR1").Copy()
$ws.Range("AS1").PasteSpecial(-4122)
$ws.Range("AS1").Value = "BOUNDARY"

# --- Row 2: epsilon values 0.01 .. 0.20, styled like the other epsilon blocks ---
$ws.Range("AK2:AR2").Copy()
$ws.Range("AS2").PasteSpecial(-4122)
$ws.Range("AK2:AR2").Copy()
$ws.Range("AS2").PasteSpecial(-4163)

# --- Rows 4-12 (except row 3): new numeric data for the BOUNDARY attack block ---
$ws.Range("AS4").Value = 3.365715980529785
$ws.Range("AT4").Value = 3.425267219543457
$ws.Range("AU4").Value = 3.436357021331787
$ws.Range("AV4").Value = 3.506391763687134
$ws.Range("AW4").Value = 3.66422963142395
$ws.Range("AX4").Value = 3.871299743652344
$ws.Range("AY4").Value = 4.294879913330078
$ws.Range("AZ4").Value = 6.093916416168213
$ws.Range("AS5").Value = 4.391582561076026
$ws.Range("AT5").Value = 4.453119003559559
$ws.Range("AU5").Value = 4.482705898885552
$ws.Range("AV5").Value = 4.57630504067609
$ws.Range("AW5").Value = 4.677706647882333
$ws.Range("AX5").Value = 5.035980654302863
$ws.Range("AY5").Value = 5.408472387155367
$ws.Range("AZ5").Value = 7.571858566083719
$ws.Range("AS6").Value = 0.9996684193611145
$ws.Range("AT6").Value = 0.9996588230133057
$ws.Range("AU6").Value = 0.9996528029441833
$ws.Range("AV6").Value = 0.9996389746665955
$ws.Range("AW6").Value = 0.9996217489242554
$ws.Range("AX6").Value = 0.9995622634887695
$ws.Range("AY6").Value = 0.9994813203811646
$ws.Range("AZ6").Value = 0.998970627784729
$ws.Range("AS7").Value = 3.456785678863525
$ws.Range("AT7").Value = 3.548597812652588
$ws.Range("AU7").Value = 3.804541110992432
$ws.Range("AV7").Value = 4.195152759552002
$ws.Range("AW7").Value = 4.487371444702148
$ws.Range("AX7").Value = 5.448037624359131
$ws.Range("AY7").Value = 6.861483573913574
$ws.Range("AZ7").Value = 11.56615543365479
$ws.Range("AS8").Value = 4.342254381189789
$ws.Range("AT8").Value = 4.439331321539029
$ws.Range("AU8").Value = 4.785290176697262
$ws.Range("AV8").Value = 5.277365435281131
$ws.Range("AW8").Value = 5.639112444091223
$ws.Range("AX8").Value = 6.805521619370461
$ws.Range("AY8").Value = 8.545227003480694
$ws.Range("AZ8").Value = 14.57730227782969
$ws.Range("AS9").Value = 0.9997920989990234
$ws.Range("AT9").Value = 0.9997710585594177
$ws.Range("AU9").Value = 0.999718189239502
$ws.Range("AV9").Value = 0.9996262788772583
$ws.Range("AW9").Value = 0.9995540976524353
$ws.Range("AX9").Value = 0.9992905259132385
$ws.Range("AY9").Value = 0.9988683462142944
$ws.Range("AZ9").Value = 0.9963282942771912
$ws.Range("AS10").Value = 3.063418388366699
$ws.Range("AT10").Value = 3.161049604415894
$ws.Range("AU10").Value = 3.251139879226685
$ws.Range("AV10").Value = 3.56424355506897
$ws.Range("AW10").Value = 3.829797983169556
$ws.Range("AX10").Value = 4.520462512969971
$ws.Range("AY10").Value = 5.508560180664062
$ws.Range("AZ10").Value = 8.935407638549805
$ws.Range("AS11").Value = 3.99095669341003
$ws.Range("AT11").Value = 4.066556343291467
$ws.Range("AU11").Value = 4.205771368265758
$ws.Range("AV11").Value = 4.511735871517837
$ws.Range("AW11").Value = 4.897603379191018
$ws.Range("AX11").Value = 5.702020240279461
$ws.Range("AY11").Value = 6.943561002698353
$ws.Range("AZ11").Value = 11.18518656248075
$ws.Range("AS12").Value = 0.9997690916061401
$ws.Range("AT12").Value = 0.9997552037239075
$ws.Range("AU12").Value = 0.9997372627258301
$ws.Range("AV12").Value = 0.9996918439865112
$ws.Range("AW12").Value = 0.9996218681335449
$ws.Range("AX12").Value = 0.9994951486587524
$ws.Range("AY12").Value = 0.999189555644989
$ws.Range("AZ12").Value = 0.9978799819946289
